# Auto-generated Excel COM-interop script
# Applies scheduled-runner market data refresh to Sheets/Ifrit_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1564.2858
$ws.Range("J17").Value = 1564.2858
$ws.Range("L17").Value = 4692.857400000001
$ws.Range("N17").Value = -5028.857400000001
$ws.Range("H64").Value = 4622.5
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 3993.3333
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 3993.3333
$ws.Range("M64").Value = -4752
$ws.Range("N64").Value = -4489.3333
$ws.Range("H67").Value = 4622.5
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 3993.3333
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 3993.3333
$ws.Range("M67").Value = -4142
$ws.Range("N67").Value = -5709.3333
$ws.Range("H74").Value = 3469.4167
$ws.Range("I74").Value = 3388.111
$ws.Range("K74").Value = 3388.111
$ws.Range("M74").Value = -2452.111
$ws.Range("H77").Value = 3469.4167
$ws.Range("I77").Value = 3388.111
$ws.Range("K77").Value = 16940.555
$ws.Range("M77").Value = -12260.555
$ws.Range("H80").Value = 387.86957
$ws.Range("I80").Value = 527.1818
$ws.Range("J80").Value = 260.16666
$ws.Range("K80").Value = 1581.5454
$ws.Range("L80").Value = 780.4999799999999
$ws.Range("M80").Value = -583.5454
$ws.Range("N80").Value = -2776.49998
$ws.Range("H83").Value = 387.86957
$ws.Range("I83").Value = 527.1818
$ws.Range("J83").Value = 260.16666
$ws.Range("K83").Value = 4744.6362
$ws.Range("L83").Value = 2341.49994
$ws.Range("M83").Value = 247.3638000000001
$ws.Range("N83").Value = -12325.49994
$ws.Range("H138").Value = 2153.69
$ws.Range("I138").Value = 931.0857
$ws.Range("J138").Value = 2812.0154
$ws.Range("K138").Value = 2793.2571
$ws.Range("L138").Value = 8436.046200000001
$ws.Range("M138").Value = 2346.7429
$ws.Range("N138").Value = -18716.0462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2446.97
$ws.Range("I32").Value = 2421.182
$ws.Range("K32").Value = 2421.182
$ws.Range("M32").Value = -2134.182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 10000
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 10000
$ws.Range("N46").Value = -10596
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 866.6667
$ws.Range("J13").Value = 300
$ws.Range("L13").Value = 300
$ws.Range("N13").Value = -578
$ws.Range("H31").Value = 1496.1702
$ws.Range("I31").Value = 1012.5128
$ws.Range("J31").Value = 3854
$ws.Range("K31").Value = 1012.5128
$ws.Range("L31").Value = 3854
$ws.Range("M31").Value = -717.5128
$ws.Range("N31").Value = -4444
$ws.Range("H34").Value = 1496.1702
$ws.Range("I34").Value = 1012.5128
$ws.Range("J34").Value = 3854
$ws.Range("K34").Value = 1012.5128
$ws.Range("L34").Value = 3854
$ws.Range("M34").Value = -810.5128
$ws.Range("N34").Value = -4258
$ws.Range("H99").Value = 3117.4546
$ws.Range("I99").Value = 3032.4443
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 3032.4443
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -1534.4443
$ws.Range("N99").Value = -6496
$ws.Range("H126").Value = 3117.4546
$ws.Range("I126").Value = 3032.4443
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 9097.332900000001
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -6627.332900000001
$ws.Range("N126").Value = -15440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2043
$ws.Range("I22").Value = 990
$ws.Range("J22").Value = 2130.75
$ws.Range("K22").Value = 2970
$ws.Range("L22").Value = 6392.25
$ws.Range("M22").Value = -2801
$ws.Range("N22").Value = -6730.25
$ws.Range("H26").Value = 2400
$ws.Range("I26").Value = 100
$ws.Range("J26").Value = 3166.6667
$ws.Range("K26").Value = 300
$ws.Range("L26").Value = 9500.000100000001
$ws.Range("M26").Value = -12
$ws.Range("N26").Value = -10076.0001
$ws.Range("H27").Value = 2043
$ws.Range("I27").Value = 990
$ws.Range("J27").Value = 2130.75
$ws.Range("K27").Value = 2970
$ws.Range("L27").Value = 6392.25
$ws.Range("M27").Value = -2868
$ws.Range("N27").Value = -6596.25
$ws.Range("H69").Value = 1198.2
$ws.Range("I69").Value = 773
$ws.Range("J69").Value = 1481.6666
$ws.Range("K69").Value = 2319
$ws.Range("L69").Value = 4444.9998
$ws.Range("M69").Value = -1508
$ws.Range("N69").Value = -6066.9998
$ws.Range("H72").Value = 1198.2
$ws.Range("I72").Value = 773
$ws.Range("J72").Value = 1481.6666
$ws.Range("K72").Value = 6957
$ws.Range("L72").Value = 13334.9994
$ws.Range("M72").Value = -2901
$ws.Range("N72").Value = -21446.9994
$ws.Range("H75").Value = 972.5
$ws.Range("J75").Value = 1068.5714
$ws.Range("L75").Value = 3205.7142
$ws.Range("N75").Value = -5201.7142
$ws.Range("H78").Value = 972.5
$ws.Range("J78").Value = 1068.5714
$ws.Range("L78").Value = 9617.142600000001
$ws.Range("N78").Value = -19601.1426
$ws.Range("H103").Value = 2834504.2
$ws.Range("I103").Value = 6800465
$ws.Range("J103").Value = 1675.1428
$ws.Range("K103").Value = 20401395
$ws.Range("L103").Value = 5025.428400000001
$ws.Range("M103").Value = -20400516
$ws.Range("N103").Value = -6783.428400000001
$ws.Range("H107").Value = 105578.266
$ws.Range("I107").Value = 111414.78
$ws.Range("J107").Value = 100325.4
$ws.Range("K107").Value = 334244.34
$ws.Range("L107").Value = 300976.2
$ws.Range("M107").Value = -332324.34
$ws.Range("N107").Value = -304816.2
$ws.Range("H114").Value = 1201.2069
$ws.Range("I114").Value = 634.5714
$ws.Range("J114").Value = 1730.0667
$ws.Range("K114").Value = 1903.7142
$ws.Range("L114").Value = 5190.2001
$ws.Range("M114").Value = 1350.2858
$ws.Range("N114").Value = -11698.2001
$ws.Range("H121").Value = 16667422
$ws.Range("I121").Value = 434.9
$ws.Range("J121").Value = 25000916
$ws.Range("K121").Value = 1304.7
$ws.Range("L121").Value = 75002748
$ws.Range("M121").Value = 5.300000000000182
$ws.Range("N121").Value = -75005368

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 4997.5
$ws.Range("I44").Value = 2495
$ws.Range("J44").Value = 7500
$ws.Range("K44").Value = 2495
$ws.Range("L44").Value = 7500
$ws.Range("M44").Value = -1899
$ws.Range("N44").Value = -8692
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 4790.8335
$ws.Range("I4").Value = 2490
$ws.Range("K4").Value = 2490
$ws.Range("M4").Value = -2377
$ws.Range("H28").Value = 4790.8335
$ws.Range("I28").Value = 2490
$ws.Range("K28").Value = 2490
$ws.Range("M28").Value = -2258
$ws.Range("H37").Value = 4790.8335
$ws.Range("I37").Value = 2490
$ws.Range("K37").Value = 2490
$ws.Range("M37").Value = -2383
$ws.Range("H40").Value = 1914.8572
$ws.Range("I40").Value = 1621.1818
$ws.Range("J40").Value = 2991.6667
$ws.Range("K40").Value = 1621.1818
$ws.Range("L40").Value = 2991.6667
$ws.Range("M40").Value = -1485.1818
$ws.Range("N40").Value = -3263.6667
$ws.Range("H100").Value = 2733.818
$ws.Range("J100").Value = 4160.5557
$ws.Range("L100").Value = 4160.5557
$ws.Range("N100").Value = -5242.5557
$ws.Range("H132").Value = 3841.5667
$ws.Range("I132").Value = 3799.926
$ws.Range("K132").Value = 11399.778
$ws.Range("M132").Value = -8869.778
$ws.Range("H136").Value = 1495.5
$ws.Range("I136").Value = 971.92
$ws.Range("J136").Value = 2685.4546
$ws.Range("K136").Value = 2915.76
$ws.Range("L136").Value = 8056.3638
$ws.Range("M136").Value = -365.7599999999998
$ws.Range("N136").Value = -13156.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 53885.8
$ws.Range("J138").Value = 53885.8
$ws.Range("L138").Value = 53885.8
$ws.Range("N138").Value = -64165.8
